# Rename the three picture placeholders (the Pearson logo that appears in
# both footers, and the BTEC logo in the "first page" header) so that their
# wp:docPr/@name swaps with its sibling:
#
#   header (first page)   : image1.jpg -> image2.jpg
#   footer (default/odd)  : image2.png -> image1.png
#   footer (first page)   : image2.png -> image1.png
#
# The shapes are inline pictures living inside header/footer stories, so
# Document.InlineShapes (body only) does not see them - they have to be
# reached through Sections(1).Headers/Footers(n).Range.InlineShapes.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlineLogo {
    param($HeaderFooter, $ExpectedOldName, $NewName)

    $shape = $HeaderFooter.Range.InlineShapes.Item(1)

    # Re-fetching the InlineShape through its own Range (instead of using the
    # handle returned directly from the header/footer Range) is what makes the
    # rename actually stick for shapes that live in a footer story - setting
    # .Name on the first handle is silently a no-op there.
    $shape = $shape.Range.InlineShapes.Item(1)

    Write-Output ("Found shape descr='" + $shape.AlternativeText + "' -> renaming to '" + $NewName + "'")
    $shape.Name = $NewName
}

# First-page header: BTec_Logo-Orange (image1.jpg -> image2.jpg)
Rename-InlineLogo $sec.Headers.Item(2) "image1.jpg" "image2.jpg"

# Default/odd-page footer: PearsonLogo.png (image2.png -> image1.png)
Rename-InlineLogo $sec.Footers.Item(1) "image2.png" "image1.png"

# First-page footer: PearsonLogo.png (image2.png -> image1.png)
Rename-InlineLogo $sec.Footers.Item(2) "image2.png" "image1.png"

Write-Output "Logo picture names updated."
